$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,3
$data[0,0] = 300.3355232157551
$data[0,1] = 1.957258766345068
$data[0,2] = -0.992936891014108
$data[1,0] = 302.1200307138862
$data[1,1] = 1.980089354100932
$data[1,2] = -0.9994303881169618
$data[2,0] = 294.9806523689379
$data[2,1] = 2.1248770935966332
$data[2,2] = -1.005548259112435
$data[3,0] = 289.15970945105363
$data[3,1] = 2.059263463833459
$data[3,2] = -0.90551621756431
$data[4,0] = 302.2546139555187
$data[4,1] = 1.9610693809251658
$data[4,2] = -0.9916869934327884
$data[5,0] = 297.9454685279235
$data[5,1] = 1.8606134021166745
$data[5,2] = -0.8810421064126592
$data[6,0] = 303.1344483809093
$data[6,1] = 1.9320072656977718
$data[6,2] = -0.8988336229418873
$data[7,0] = 310.29549349639655
$data[7,1] = 1.9039744372002845
$data[7,2] = -1.109354809402632
$data[8,0] = 294.6867198194745
$data[8,1] = 2.0545402293367534
$data[8,2] = -1.053924603653631
$data[9,0] = 292.3151801403916
$data[9,1] = 2.030782533246721
$data[9,2] = -0.9499285749719409
$data[10,0] = 310.92379620920235
$data[10,1] = 1.9523161365713426
$data[10,2] = -1.080673124955735
$data[11,0] = 292.7880215583621
$data[11,1] = 1.9769105668579292
$data[11,2] = -0.8352055958792949
$data[12,0] = 310.96020431255886
$data[12,1] = 1.8850967138003891
$data[12,2] = -1.0805338955682469
$data[13,0] = 291.8318779832073
$data[13,1] = 2.168460844260386
$data[13,2] = -1.0157390263736483
$data[14,0] = 299.6280854987036
$data[14,1] = 1.9675995814451777
$data[14,2] = -1.0442245163717248
$data[15,0] = 296.03552154118034
$data[15,1] = 1.9872847762585253
$data[15,2] = -0.8896363166095327
$data[16,0] = 281.4968822759294
$data[16,1] = 2.262254436321159
$data[16,2] = -0.9670880146686011
$data[17,0] = 289.7510422475379
$data[17,1] = 2.1320272733355914
$data[17,2] = -0.9176287144317584
$data[18,0] = 300.9899011235117
$data[18,1] = 2.0407253467456536
$data[18,2] = -0.9470525681715379
$data[19,0] = 314.5510878791318
$data[19,1] = 1.8368188990294851
$data[19,2] = -1.0612728776735716
$data[20,0] = 309.37266044769
$data[20,1] = 1.962303121056252
$data[20,2] = -1.2230256489273788
$data[21,0] = 290.8685550266593
$data[21,1] = 1.7328592134996328
$data[21,2] = -0.7415550546033299
$data[22,0] = 307.4661665965745
$data[22,1] = 1.8916815140183176
$data[22,2] = -1.0464894413129464
$data[23,0] = 286.130226637444
$data[23,1] = 2.053339628956805
$data[23,2] = -0.847251885582815
$data[24,0] = 291.56081428278924
$data[24,1] = 2.1026098217250664
$data[24,2] = -0.9187012329173585
$data[25,0] = 280.74453440569374
$data[25,1] = 2.142512501556268
$data[25,2] = -0.7445160704077327
$data[26,0] = 298.5983602861604
$data[26,1] = 2.067537883277829
$data[26,2] = -0.949949721865206
$data[27,0] = 306.4441611523221
$data[27,1] = 1.997686802816628
$data[27,2] = -1.2795351469796803
$data[28,0] = 313.6710454585673
$data[28,1] = 1.8970643097674358
$data[28,2] = -1.1744243833540473
$data[29,0] = 306.4257822135574
$data[29,1] = 1.6451993185596387
$data[29,2] = -1.061695474348634
$data[30,0] = 306.77379100771645
$data[30,1] = 1.8504152560696963
$data[30,2] = -0.9517126569887219
$data[31,0] = 299.308409746554
$data[31,1] = 1.8690669403358098
$data[31,2] = -0.8700582158107949
$data[32,0] = 284.8599381559061
$data[32,1] = 2.0868168454552887
$data[32,2] = -0.6939105074909117
$data[33,0] = 289.61163173933676
$data[33,1] = 2.202873862545679
$data[33,2] = -0.9017881676650067
$data[34,0] = 317.39431312219824
$data[34,1] = 2.0550300801755474
$data[34,2] = -1.18847985427789
$data[35,0] = 300.7989692939087
$data[35,1] = 1.9462252028006335
$data[35,2] = -0.8341201110147657
$data[36,0] = 298.3239323386273
$data[36,1] = 2.031401409946047
$data[36,2] = -0.9036463576448129
$data[37,0] = 302.71891016388304
$data[37,1] = 1.9644733930312515
$data[37,2] = -1.0951252678531198
$data[38,0] = 283.1878730597422
$data[38,1] = 2.0876519191470213
$data[38,2] = -0.8811489338448948
$data[39,0] = 311.87880288575917
$data[39,1] = 1.806598127398824
$data[39,2] = -1.0793389056618605
$data[40,0] = 302.9565813451401
$data[40,1] = 1.993150432715906
$data[40,2] = -1.0531365733667242
$data[41,0] = 296.1524203127572
$data[41,1] = 1.959602768288211
$data[41,2] = -0.95255380053481
$data[42,0] = 304.1785610371654
$data[42,1] = 2.0147758924845425
$data[42,2] = -1.02093745746713
$data[43,0] = 299.3055371877944
$data[43,1] = 2.028598707105818
$data[43,2] = -1.052427054073405
$data[44,0] = 306.8767724289121
$data[44,1] = 1.8562698350065303
$data[44,2] = -1.0336616931180485
$data[45,0] = 313.4969054504496
$data[45,1] = 1.9906804906896645
$data[45,2] = -1.1516786159257915
$data[46,0] = 327.1735481510344
$data[46,1] = 1.8260559478626412
$data[46,2] = -1.1976155147061642
$data[47,0] = 291.7469532430533
$data[47,1] = 2.326455427492738
$data[47,2] = -1.1265704347145413
$data[48,0] = 299.3478342051676
$data[48,1] = 2.1513291173058073
$data[48,2] = -1.1291722852937272
$data[49,0] = 282.72988090050416
$data[49,1] = 2.0256078860126054
$data[49,2] = -0.7394940549859654

$ws.Range("A2:C51").Value = $data
